$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue "D2" "96.974.01"
Set-TextValue "E2" "  +3.92%  "
Set-TextValue "D3" "3.142.58"
Set-TextValue "E3" "  +0.46%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "242.38"
Set-TextValue "E5" "  +1.63%  "
Set-TextValue "D6" "613.05"
Set-TextValue "E7" "  +2.03%  "
Set-TextValue "E8" "  -2.06%  "
Set-TextValue "E9" "  +0.04%  "
Set-TextValue "D10" "3.139.73"
Set-TextValue "E10" "  +0.46%  "
Set-TextValue "E11" "  -4.04%  "
Set-TextValue "E12" "  +0.18%  "
Set-TextValue "D13" "96.729.47"
Set-TextValue "E13" "  +3.98%  "
Set-TextValue "B14" "ShibaInu"
Set-TextValue "C14" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D14" "0.0000241"
Set-TextValue "E14" "  -2.15%  "
Set-TextValue "B15" "Toncoin"
Set-TextValue "C15" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D15" "5.57"
Set-TextValue "E15" "  +2.06%  "
Set-TextValue "D16" "34.22"
Set-TextValue "E16" "  -1.94%  "
Set-TextValue "D17" "3.723.70"
Set-TextValue "E17" "  +0.33%  "
Set-TextValue "D18" "3.157.38"
Set-TextValue "E18" "  +0.71%  "
Set-TextValue "E19" "  -6.24%  "
Set-TextValue "D20" "519.27"
Set-TextValue "E20" "  +17.22%  "
Set-TextValue "D21" "14.65"
Set-TextValue "E21" "  -0.13%  "
Set-TextValue "D22" "5.68"
Set-TextValue "E22" "  -4.69%  "
Set-TextValue "E23" "  -5.42%  "
Set-TextValue "D24" "8.83"
Set-TextValue "E24" "  -3.70%  "
Set-TextValue "D25" "5.51"
Set-TextValue "E25" "  -2.49%  "
Set-TextValue "D26" "88.83"
Set-TextValue "E26" "  +3.31%  "
Set-TextValue "D27" "11.67"
Set-TextValue "E27" "  -6.63%  "
Set-TextValue "D28" "3.302.56"
Set-TextValue "E28" "  +0.28%  "
Set-TextValue "E29" "  -0.11%  "
Set-TextValue "E30" "  +2.53%  "
Set-TextValue "E31" "  -2.14%  "
Set-TextValue "D32" "0.124"
Set-TextValue "E32" "  +1.07%  "
Set-TextValue "E33" "  -2.98%  "
Set-TextValue "E34" "  -1.72%  "
Set-TextValue "D35" "26.71"
Set-TextValue "E35" "  +3.08%  "
Set-TextValue "E36" "  -4.49%  "
Set-TextValue "D37" "7.38"
Set-TextValue "E37" "  -9.45%  "
Set-TextValue "E38" "  -1.06%  "
Set-TextValue "D39" "24.22"
Set-TextValue "E39" "  +0.95%  "
Set-TextValue "D40" "480.57"
Set-TextValue "E40" "  +1.02%  "
Set-TextValue "D41" "0.438"
Set-TextValue "E41" "  +1.16%  "
Set-TextValue "D42" "1.24"
Set-TextValue "E42" "  -4.78%  "
Set-TextValue "E43" "  -10.14%  "
Set-TextValue "E44" "  +0.00%  "
Set-TextValue "E45" "  -4.89%  "
Set-TextValue "D46" "161.30"
Set-TextValue "E46" "  +1.34%  "
Set-TextValue "D47" "0.709"
Set-TextValue "E47" "  +2.43%  "
Set-TextValue "D48" "1.93"
Set-TextValue "E48" "  +4.30%  "
Set-TextValue "E49" "  +2.11%  "
Set-TextValue "D50" "44.31"
Set-TextValue "E50" "  +0.44%  "
Set-TextValue "E51" "  -0.02%  "
